$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 4129.0645
$ws.Range("I76").Value = 3333
$ws.Range("J76").Value = 4454.727
$ws.Range("K76").Value = 3333
$ws.Range("L76").Value = 4454.727
$ws.Range("M76").Value = -3018
$ws.Range("N76").Value = -5084.727
# Row 79
$ws.Range("H79").Value = 4129.0645
$ws.Range("I79").Value = 3333
$ws.Range("J79").Value = 4454.727
$ws.Range("K79").Value = 3333
$ws.Range("L79").Value = 4454.727
$ws.Range("M79").Value = -2241
$ws.Range("N79").Value = -6638.727
# Row 132
$ws.Range("H132").Value = 39011240
$ws.Range("I132").Value = 45106580
$ws.Range("J132").Value = 1080.4
$ws.Range("K132").Value = 135319740
$ws.Range("L132").Value = 3241.2
$ws.Range("M132").Value = -135317210
$ws.Range("N132").Value = -8301.2
# Row 137
$ws.Range("H137").Value = 20732.77
$ws.Range("I137").Value = 1151.3
$ws.Range("J137").Value = 32971.188
$ws.Range("K137").Value = 3453.9
$ws.Range("L137").Value = 98913.56400000001
$ws.Range("M137").Value = -903.8999999999996
$ws.Range("N137").Value = -104013.564

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 5025.1113
$ws.Range("I45").Value = 9103
$ws.Range("J45").Value = 1762.8
$ws.Range("K45").Value = 9103
$ws.Range("L45").Value = 1762.8
$ws.Range("M45").Value = -8726
$ws.Range("N45").Value = -2516.8
# Row 52
$ws.Range("H52").Value = 46200
$ws.Range("J52").Value = 46200
$ws.Range("L52").Value = 46200
$ws.Range("N52").Value = -46836
# Row 61
$ws.Range("H61").Value = 4065.6304
$ws.Range("I61").Value = 4891.4243
$ws.Range("J61").Value = 1969.3846
$ws.Range("K61").Value = 4891.4243
$ws.Range("L61").Value = 1969.3846
$ws.Range("M61").Value = -4679.4243
$ws.Range("N61").Value = -2393.3846
# Row 102
$ws.Range("H102").Value = 500001500
$ws.Range("J102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("N102").Value = -6244
# Row 122
$ws.Range("H122").Value = 85715600
$ws.Range("I122").Value = 92308870
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 276926610
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -276924160
$ws.Range("N122").Value = -13900
# Row 136
$ws.Range("H136").Value = 4065.6304
$ws.Range("I136").Value = 4891.4243
$ws.Range("J136").Value = 1969.3846
$ws.Range("K136").Value = 14674.2729
$ws.Range("L136").Value = 5908.1538
$ws.Range("M136").Value = -12124.2729
$ws.Range("N136").Value = -11008.1538

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 103
$ws.Range("H103").Value = 27210.143
$ws.Range("J103").Value = 27210.143
$ws.Range("L103").Value = 27210.143
$ws.Range("N103").Value = -29554.143
# Row 105
$ws.Range("H105").Value = 2085.2666
$ws.Range("I105").Value = 2036.9231
$ws.Range("J105").Value = 2399.5
$ws.Range("K105").Value = 2036.9231
$ws.Range("L105").Value = 2399.5
$ws.Range("M105").Value = -289.9231
$ws.Range("N105").Value = -5893.5
# Row 107
$ws.Range("H107").Value = 14511.904
$ws.Range("I107").Value = 9425
$ws.Range("J107").Value = 24685.715
$ws.Range("K107").Value = 9425
$ws.Range("L107").Value = 24685.715
$ws.Range("M107").Value = -7505
$ws.Range("N107").Value = -28525.715

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7170.4287
$ws.Range("I31").Value = 3328.7856
$ws.Range("J31").Value = 11780.4
$ws.Range("K31").Value = 3328.7856
$ws.Range("L31").Value = 11780.4
$ws.Range("M31").Value = -3033.7856
$ws.Range("N31").Value = -12370.4
# Row 34
$ws.Range("H34").Value = 7170.4287
$ws.Range("I34").Value = 3328.7856
$ws.Range("J34").Value = 11780.4
$ws.Range("K34").Value = 3328.7856
$ws.Range("L34").Value = 11780.4
$ws.Range("M34").Value = -3126.7856
$ws.Range("N34").Value = -12184.4
# Row 58
$ws.Range("H58").Value = 3182775.8
$ws.Range("I58").Value = 5714973
$ws.Range("J58").Value = 17528.9
$ws.Range("K58").Value = 5714973
$ws.Range("L58").Value = 17528.9
$ws.Range("M58").Value = -5714770
$ws.Range("N58").Value = -17934.9
# Row 94
$ws.Range("H94").Value = 1280.3529
$ws.Range("I94").Value = 818.4
$ws.Range("J94").Value = 1472.8334
$ws.Range("K94").Value = 818.4
$ws.Range("L94").Value = 1472.8334
$ws.Range("M94").Value = -367.4
$ws.Range("N94").Value = -2374.8334
# Row 99
$ws.Range("H99").Value = 2804.8333
$ws.Range("I99").Value = 2295.0908
$ws.Range("J99").Value = 3605.8572
$ws.Range("K99").Value = 2295.0908
$ws.Range("L99").Value = 3605.8572
$ws.Range("M99").Value = -797.0908
$ws.Range("N99").Value = -6601.8572
# Row 126
$ws.Range("H126").Value = 2804.8333
$ws.Range("I126").Value = 2295.0908
$ws.Range("J126").Value = 3605.8572
$ws.Range("K126").Value = 6885.2724
$ws.Range("L126").Value = 10817.5716
$ws.Range("M126").Value = -4415.2724
$ws.Range("N126").Value = -15757.5716
# Row 132
$ws.Range("H132").Value = 6806718
$ws.Range("I132").Value = 10102611
$ws.Range("J132").Value = 8938.6875
$ws.Range("K132").Value = 30307833
$ws.Range("L132").Value = 26816.0625
$ws.Range("M132").Value = -30305303
$ws.Range("N132").Value = -31876.0625
# Row 136
$ws.Range("H136").Value = 3182775.8
$ws.Range("I136").Value = 5714973
$ws.Range("J136").Value = 17528.9
$ws.Range("K136").Value = 17144919
$ws.Range("L136").Value = 52586.7
$ws.Range("M136").Value = -17142369
$ws.Range("N136").Value = -57686.7

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2337.5925
$ws.Range("I68").Value = 645.1667
$ws.Range("J68").Value = 4160.205
$ws.Range("K68").Value = 1935.5001
$ws.Range("L68").Value = 12480.615
$ws.Range("M68").Value = -1124.5001
$ws.Range("N68").Value = -14102.615
# Row 71
$ws.Range("H71").Value = 2337.5925
$ws.Range("I71").Value = 645.1667
$ws.Range("J71").Value = 4160.205
$ws.Range("K71").Value = 5806.5003
$ws.Range("L71").Value = 37441.845
$ws.Range("M71").Value = -1750.5003
$ws.Range("N71").Value = -45553.845
# Row 107
$ws.Range("H107").Value = 816.15
$ws.Range("J107").Value = 2716.0908
$ws.Range("L107").Value = 8148.2724
$ws.Range("N107").Value = -11988.2724
# Row 113
$ws.Range("H113").Value = 4167158.2
$ws.Range("I113").Value = 486.66666
$ws.Range("J113").Value = 6667161
$ws.Range("K113").Value = 1459.99998
$ws.Range("L113").Value = 20001483
$ws.Range("M113").Value = 710.00002
$ws.Range("N113").Value = -20005823
# Row 118
$ws.Range("H118").Value = 125001350
$ws.Range("I118").Value = 333333600
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 1000000800
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = -999999557
$ws.Range("N118").Value = -8486
# Row 125
$ws.Range("H125").Value = 6261.1113
$ws.Range("I125").Value = 4283.3335
$ws.Range("J125").Value = 7250
$ws.Range("K125").Value = 12850.0005
$ws.Range("L125").Value = 21750
$ws.Range("M125").Value = -7930.000499999998
$ws.Range("N125").Value = -31590
# Row 131
$ws.Range("H131").Value = 21669294
$ws.Range("I131").Value = 90909450
$ws.Range("J131").Value = 2140019.2
$ws.Range("K131").Value = 272728350
$ws.Range("L131").Value = 6420057.600000001
$ws.Range("M131").Value = -272723310
$ws.Range("N131").Value = -6430137.600000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 111112160
$ws.Range("I97").Value = 66667660
$ws.Range("J97").Value = 333334620
$ws.Range("K97").Value = 66667660
$ws.Range("L97").Value = 333334620
$ws.Range("M97").Value = -66667164
$ws.Range("N97").Value = -333335612
# Row 100
$ws.Range("H100").Value = 35355
$ws.Range("J100").Value = 35355
$ws.Range("L100").Value = 35355
$ws.Range("N100").Value = -37519
# Row 130
$ws.Range("H130").Value = 50780
$ws.Range("J130").Value = 50780
$ws.Range("L130").Value = 50780
$ws.Range("N130").Value = -60820

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1517.2727
$ws.Range("I40").Value = 1319
$ws.Range("K40").Value = 1319
$ws.Range("M40").Value = -1183
# Row 122
$ws.Range("H122").Value = 31253726
$ws.Range("I122").Value = 1152
$ws.Range("J122").Value = 41671250
$ws.Range("K122").Value = 3456
$ws.Range("L122").Value = 125013750
$ws.Range("M122").Value = -1006
$ws.Range("N122").Value = -125018650
# Row 130
$ws.Range("H130").Value = 34429
$ws.Range("J130").Value = 34429
$ws.Range("L130").Value = 34429
$ws.Range("N130").Value = -44469
# Row 132
$ws.Range("H132").Value = 4123018.8
$ws.Range("I132").Value = 4117728.5
$ws.Range("J132").Value = 4136004.2
$ws.Range("K132").Value = 12353185.5
$ws.Range("L132").Value = 12408012.6
$ws.Range("M132").Value = -12350655.5
$ws.Range("N132").Value = -12413072.6
# Row 135
$ws.Range("H135").Value = 39429
$ws.Range("J135").Value = 39429
$ws.Range("L135").Value = 39429
$ws.Range("N135").Value = -49569

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 8481.154
$ws.Range("I122").Value = 9841.363
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 29524.089
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -27074.089
$ws.Range("N122").Value = -7900
